$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H15").Value = 184465.45
$ws.Range("I15").Value = 184465.45
$ws.Range("K15").Value = 553396.3500000001
$ws.Range("M15").Value = -553227.3500000001

$ws.Range("H17").Value = 1613638
$ws.Range("J17").Value = 1613638
$ws.Range("L17").Value = 4840914
$ws.Range("N17").Value = -4841250

$ws.Range("H86").Value = 19963.5
$ws.Range("I86").Value = 5901
$ws.Range("J86").Value = 58133.145
$ws.Range("K86").Value = 5901
$ws.Range("L86").Value = 58133.145
$ws.Range("M86").Value = -4778
$ws.Range("N86").Value = -60379.145

$ws.Range("H87").Value = 27023.219
$ws.Range("J87").Value = 27023.219
$ws.Range("L87").Value = 27023.219
$ws.Range("N87").Value = -29519.219

$ws.Range("H89").Value = 19963.5
$ws.Range("I89").Value = 5901
$ws.Range("J89").Value = 58133.145
$ws.Range("K89").Value = 29505
$ws.Range("L89").Value = 290665.725
$ws.Range("M89").Value = -23889
$ws.Range("N89").Value = -301897.725

$ws.Range("H90").Value = 27023.219
$ws.Range("J90").Value = 27023.219
$ws.Range("L90").Value = 81069.65700000001
$ws.Range("N90").Value = -93549.65700000001

$ws.Range("H92").Value = 377.85
$ws.Range("I92").Value = 332.05554
$ws.Range("J92").Value = 790
$ws.Range("K92").Value = 332.05554
$ws.Range("L92").Value = 790
$ws.Range("M92").Value = 915.9444599999999
$ws.Range("N92").Value = -3286

$ws.Range("H100").Value = 11943017
$ws.Range("I100").Value = 16668229
$ws.Range("J100").Value = 129988.125
$ws.Range("K100").Value = 16668229
$ws.Range("L100").Value = 129988.125
$ws.Range("M100").Value = -16667688
$ws.Range("N100").Value = -131070.125

$ws.Range("H135").Value = 350.69565
$ws.Range("I135").Value = 326.97675
$ws.Range("J135").Value = 690.6667
$ws.Range("K135").Value = 2942.79075
$ws.Range("L135").Value = 6216.0003
$ws.Range("M135").Value = -407.7907499999997
$ws.Range("N135").Value = -11286.0003

$ws.Range("H138").Value = 2770.9158
$ws.Range("I138").Value = 1315.3392
$ws.Range("J138").Value = 4860.974
$ws.Range("K138").Value = 3946.0176
$ws.Range("L138").Value = 14582.922
$ws.Range("M138").Value = 1193.9824
$ws.Range("N138").Value = -24862.922

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 1739.82
$ws.Range("I32").Value = 1704.102
$ws.Range("K32").Value = 1704.102
$ws.Range("M32").Value = -1417.102

$ws.Range("H45").Value = 3000.6667
$ws.Range("I45").Value = 5012
$ws.Range("J45").Value = 1995
$ws.Range("K45").Value = 5012
$ws.Range("L45").Value = 1995
$ws.Range("M45").Value = -4635
$ws.Range("N45").Value = -2749

$ws.Range("H110").Value = 53539.2
$ws.Range("I110").Value = 127730.5
$ws.Range("J110").Value = 4078.3333
$ws.Range("K110").Value = 127730.5
$ws.Range("L110").Value = 4078.3333
$ws.Range("M110").Value = -125685.5
$ws.Range("N110").Value = -8168.3333

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 2712.8635
$ws.Range("I86").Value = 2927.9167
$ws.Range("J86").Value = 2454.8
$ws.Range("K86").Value = 2927.9167
$ws.Range("L86").Value = 2454.8
$ws.Range("M86").Value = -1804.9167
$ws.Range("N86").Value = -4700.8

$ws.Range("H89").Value = 2712.8635
$ws.Range("I89").Value = 2927.9167
$ws.Range("J89").Value = 2454.8
$ws.Range("K89").Value = 14639.5835
$ws.Range("L89").Value = 12274
$ws.Range("M89").Value = -9023.583500000001
$ws.Range("N89").Value = -23506

$ws.Range("H134").Value = 1869.75
$ws.Range("I134").Value = 1313.4073
$ws.Range("J134").Value = 3025.2307
$ws.Range("K134").Value = 3940.2219
$ws.Range("L134").Value = 9075.6921
$ws.Range("M134").Value = -1405.2219
$ws.Range("N134").Value = -14145.6921

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 1536.1777
$ws.Range("I58").Value = 1130.2354
$ws.Range("J58").Value = 2790.9092
$ws.Range("K58").Value = 1130.2354
$ws.Range("L58").Value = 2790.9092
$ws.Range("M58").Value = -927.2354
$ws.Range("N58").Value = -3196.9092

$ws.Range("H99").Value = 10013.467
$ws.Range("I99").Value = 2792
$ws.Range("J99").Value = 13624.2
$ws.Range("K99").Value = 2792
$ws.Range("L99").Value = 13624.2
$ws.Range("M99").Value = -1294
$ws.Range("N99").Value = -16620.2

$ws.Range("H107").Value = 1106.8
$ws.Range("I107").Value = 678
$ws.Range("J107").Value = 1750
$ws.Range("K107").Value = 678
$ws.Range("L107").Value = 1750
$ws.Range("M107").Value = 1242
$ws.Range("N107").Value = -5590

$ws.Range("H126").Value = 10013.467
$ws.Range("I126").Value = 2792
$ws.Range("J126").Value = 13624.2
$ws.Range("K126").Value = 8376
$ws.Range("L126").Value = 40872.60000000001
$ws.Range("M126").Value = -5906
$ws.Range("N126").Value = -45812.60000000001

$ws.Range("H132").Value = 2211.3235
$ws.Range("I132").Value = 1495.8148
$ws.Range("J132").Value = 4971.143
$ws.Range("K132").Value = 4487.4444
$ws.Range("L132").Value = 14913.429
$ws.Range("M132").Value = -1957.4444
$ws.Range("N132").Value = -19973.429

$ws.Range("H136").Value = 1536.1777
$ws.Range("I136").Value = 1130.2354
$ws.Range("J136").Value = 2790.9092
$ws.Range("K136").Value = 3390.7062
$ws.Range("L136").Value = 8372.7276
$ws.Range("M136").Value = -840.7062000000001
$ws.Range("N136").Value = -13472.7276

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H113").Value = 866.56665
$ws.Range("I113").Value = 1107.2941
$ws.Range("J113").Value = 551.7692
$ws.Range("K113").Value = 3321.8823
$ws.Range("L113").Value = 1655.3076
$ws.Range("M113").Value = -1151.8823
$ws.Range("N113").Value = -5995.3076

$ws.Range("H131").Value = 5101.7144
$ws.Range("J131").Value = 5653.2256
$ws.Range("L131").Value = 16959.6768
$ws.Range("N131").Value = -27039.6768

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 2576.16
$ws.Range("I122").Value = 2035.0769
$ws.Range("J122").Value = 3162.3333
$ws.Range("K122").Value = 6105.2307
$ws.Range("L122").Value = 9486.999899999999
$ws.Range("M122").Value = -3655.2307
$ws.Range("N122").Value = -14386.9999

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 1131.1111
$ws.Range("I22").Value = 650
$ws.Range("J22").Value = 1268.5714
$ws.Range("K22").Value = 650
$ws.Range("L22").Value = 1268.5714
$ws.Range("M22").Value = -355
$ws.Range("N22").Value = -1858.5714

$ws.Range("H27").Value = 1131.1111
$ws.Range("I27").Value = 650
$ws.Range("J27").Value = 1268.5714
$ws.Range("K27").Value = 650
$ws.Range("L27").Value = 1268.5714
$ws.Range("M27").Value = -543
$ws.Range("N27").Value = -1482.5714

$ws.Range("H82").Value = 1798.3572
$ws.Range("I82").Value = 1297
$ws.Range("J82").Value = 1998.9
$ws.Range("K82").Value = 1297
$ws.Range("L82").Value = 1998.9
$ws.Range("M82").Value = -936
$ws.Range("N82").Value = -2720.9

$ws.Range("H85").Value = 1798.3572
$ws.Range("I85").Value = 1297
$ws.Range("J85").Value = 1998.9
$ws.Range("K85").Value = 1297
$ws.Range("L85").Value = 1998.9
$ws.Range("M85").Value = -49
$ws.Range("N85").Value = -4494.9

$ws.Range("H132").Value = 4586.5747
$ws.Range("I132").Value = 4551.2153
$ws.Range("J132").Value = 4691.0454
$ws.Range("K132").Value = 13653.6459
$ws.Range("L132").Value = 14073.1362
$ws.Range("M132").Value = -11123.6459
$ws.Range("N132").Value = -19133.1362

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H37").Value = 0
$ws.Range("I37").Value = 0
$ws.Range("J37").Value = 0
$ws.Range("K37").Value = 0
$ws.Range("L37").Value = 0
$ws.Range("M37").ClearContents()
$ws.Range("N37").ClearContents()

$ws.Range("H132").Value = 1694.7273
$ws.Range("I132").Value = 1538.4468
$ws.Range("J132").Value = 2612.875
$ws.Range("K132").Value = 4615.3404
$ws.Range("L132").Value = 7838.625
$ws.Range("M132").Value = -2085.3404
$ws.Range("N132").Value = -12898.625
